$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 56490.5
$ws.Range("J28").Value = 989
$ws.Range("L28").Value = 989
$ws.Range("N28").Value = -1959
$ws.Range("H34").Value = 3197.5454
$ws.Range("I34").Value = 3197.5454
$ws.Range("K34").Value = 3197.5454
$ws.Range("M34").Value = -2994.5454
$ws.Range("H36").Value = 3197.5454
$ws.Range("I36").Value = 3197.5454
$ws.Range("K36").Value = 3197.5454
$ws.Range("M36").Value = -2482.5454
$ws.Range("H135").Value = 840.1818
$ws.Range("J135").Value = 4992.6665
$ws.Range("L135").Value = 44933.9985
$ws.Range("N135").Value = -50003.9985
$ws.Range("H137").Value = 2392.2222
$ws.Range("I137").Value = 2642.7144
$ws.Range("J137").Value = 2304.55
$ws.Range("K137").Value = 7928.1432
$ws.Range("L137").Value = 6913.650000000001
$ws.Range("M137").Value = -5378.1432
$ws.Range("N137").Value = -12013.65
$ws.Range("H138").Value = 1832.3636
$ws.Range("I138").Value = 1501.375
$ws.Range("J138").Value = 2021.5
$ws.Range("K138").Value = 4504.125
$ws.Range("L138").Value = 6064.5
$ws.Range("M138").Value = 635.875
$ws.Range("N138").Value = -16344.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1722.2222
$ws.Range("I45").Value = 1243.2667
$ws.Range("J45").Value = 4117
$ws.Range("K45").Value = 1243.2667
$ws.Range("L45").Value = 4117
$ws.Range("M45").Value = -866.2666999999999
$ws.Range("N45").Value = -4871
$ws.Range("H61").Value = 2915.5
$ws.Range("I61").Value = 1860.1666
$ws.Range("J61").Value = 4498.5
$ws.Range("K61").Value = 1860.1666
$ws.Range("L61").Value = 4498.5
$ws.Range("M61").Value = -1648.1666
$ws.Range("N61").Value = -4922.5
$ws.Range("H122").Value = 2588.4814
$ws.Range("I122").Value = 1583.875
$ws.Range("J122").Value = 4049.7273
$ws.Range("K122").Value = 4751.625
$ws.Range("L122").Value = 12149.1819
$ws.Range("M122").Value = -2301.625
$ws.Range("N122").Value = -17049.1819
$ws.Range("H132").Value = 6717
$ws.Range("I132").Value = 6707.15
$ws.Range("J132").Value = 6914
$ws.Range("K132").Value = 20121.45
$ws.Range("L132").Value = 20742
$ws.Range("M132").Value = -17591.45
$ws.Range("N132").Value = -25802
$ws.Range("H136").Value = 2915.5
$ws.Range("I136").Value = 1860.1666
$ws.Range("J136").Value = 4498.5
$ws.Range("K136").Value = 5580.4998
$ws.Range("L136").Value = 13495.5
$ws.Range("M136").Value = -3030.4998
$ws.Range("N136").Value = -18595.5
$ws.Range("H137").Value = 50517.332
$ws.Range("I137").Value = 50000
$ws.Range("K137").Value = 50000
$ws.Range("M137").Value = -44900

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 5725
$ws.Range("J18").Value = 5725
$ws.Range("L18").Value = 5725
$ws.Range("N18").Value = -6783
$ws.Range("H20").Value = 6780.7856
$ws.Range("I20").Value = 8407.906000000001
$ws.Range("K20").Value = 8407.906000000001
$ws.Range("M20").Value = -8160.906000000001
$ws.Range("H134").Value = 3482.6667
$ws.Range("I134").Value = 3074.3125
$ws.Range("K134").Value = 9222.9375
$ws.Range("M134").Value = -6687.9375

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2538.9167
$ws.Range("I58").Value = 1843
$ws.Range("J58").Value = 2886.875
$ws.Range("K58").Value = 1843
$ws.Range("L58").Value = 2886.875
$ws.Range("M58").Value = -1640
$ws.Range("N58").Value = -3292.875
$ws.Range("H107").Value = 488.04
$ws.Range("I107").Value = 490.8421
$ws.Range("K107").Value = 490.8421
$ws.Range("M107").Value = 1429.1579
$ws.Range("H125").Value = 36666.332
$ws.Range("J125").Value = 36666.332
$ws.Range("L125").Value = 36666.332
$ws.Range("N125").Value = -41586.332
$ws.Range("H132").Value = 3933.1428
$ws.Range("H136").Value = 2538.9167
$ws.Range("I136").Value = 1843
$ws.Range("J136").Value = 2886.875
$ws.Range("K136").Value = 5529
$ws.Range("L136").Value = 8660.625
$ws.Range("M136").Value = -2979
$ws.Range("N136").Value = -13760.625
$ws.Range("H141").Value = 192046.11
$ws.Range("J141").Value = 192046.11
$ws.Range("L141").Value = 192046.11
$ws.Range("N141").Value = -202406.11

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 767.6667
$ws.Range("J92").Value = 801.5
$ws.Range("L92").Value = 2404.5
$ws.Range("N92").Value = -4900.5
$ws.Range("H97").Value = 430.75
$ws.Range("I97").Value = 220
$ws.Range("K97").Value = 660
$ws.Range("M97").Value = -164
$ws.Range("H101").Value = 7000
$ws.Range("J101").Value = 7000
$ws.Range("L101").Value = 21000
$ws.Range("N101").Value = -25868
$ws.Range("H105").Value = 13249.25
$ws.Range("H121").Value = 936596.4
$ws.Range("I121").Value = 1133478.9
$ws.Range("K121").Value = 3400436.7
$ws.Range("M121").Value = -3399126.7
$ws.Range("H129").Value = 1620.7693
$ws.Range("I129").Value = 902.75
$ws.Range("J129").Value = 2769.6
$ws.Range("K129").Value = 2708.25
$ws.Range("L129").Value = 8308.799999999999
$ws.Range("M129").Value = 2291.75
$ws.Range("N129").Value = -18308.8

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H63").Value = 35000
$ws.Range("J63").Value = 35000
$ws.Range("L63").Value = 35000
$ws.Range("N63").Value = -36372
$ws.Range("H66").Value = 35000
$ws.Range("J66").Value = 35000
$ws.Range("L66").Value = 105000
$ws.Range("N66").Value = -111864
$ws.Range("H80").Value = 3523.08
$ws.Range("I80").Value = 1431.4445
$ws.Range("J80").Value = 4699.625
$ws.Range("K80").Value = 1431.4445
$ws.Range("L80").Value = 4699.625
$ws.Range("M80").Value = -433.4445000000001
$ws.Range("N80").Value = -6695.625
$ws.Range("H83").Value = 3523.08
$ws.Range("I83").Value = 1431.4445
$ws.Range("J83").Value = 4699.625
$ws.Range("K83").Value = 7157.2225
$ws.Range("L83").Value = 23498.125
$ws.Range("M83").Value = -2165.2225
$ws.Range("N83").Value = -33482.125
$ws.Range("H103").Value = 61650
$ws.Range("J103").Value = 61650
$ws.Range("L103").Value = 61650
$ws.Range("N103").Value = -63994
$ws.Range("H122").Value = 2099.36
$ws.Range("I122").Value = 1817.8
$ws.Range("K122").Value = 5453.4
$ws.Range("M122").Value = -3003.4
$ws.Range("H132").Value = 3401.4707
$ws.Range("I132").Value = 3855.5
$ws.Range("K132").Value = 11566.5
$ws.Range("M132").Value = -9036.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 749.5
$ws.Range("I16").Value = 499
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 499
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -329
$ws.Range("N16").Value = -1340
$ws.Range("H22").Value = 1256.875
$ws.Range("I22").Value = 2000
$ws.Range("K22").Value = 2000
$ws.Range("M22").Value = -1705
$ws.Range("H27").Value = 1256.875
$ws.Range("I27").Value = 2000
$ws.Range("K27").Value = 2000
$ws.Range("M27").Value = -1893
$ws.Range("H34").Value = 5000
$ws.Range("I34").Value = 5000
$ws.Range("K34").Value = 5000
$ws.Range("M34").Value = -4828

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 50090.332
$ws.Range("J69").Value = 50090.332
$ws.Range("L69").Value = 50090.332
$ws.Range("N69").Value = -51588.332
$ws.Range("H72").Value = 50090.332
$ws.Range("J72").Value = 50090.332
$ws.Range("L72").Value = 150270.996
$ws.Range("N72").Value = -157758.996
$ws.Range("H135").Value = 40704.668
$ws.Range("J135").Value = 40704.668
$ws.Range("L135").Value = 40704.668
$ws.Range("N135").Value = -50844.668
$ws.Range("H136").Value = 2467.476
$ws.Range("I136").Value = 2051.7144
$ws.Range("K136").Value = 6155.1432
$ws.Range("M136").Value = -3605.1432
